$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 159; A = 11; B = "Vega Monumental Concepción"; C = "Bíobío"; D = 44595; E = 8; F = "Fruta"; G = 100103; H = "Frutos de hueso (carozo)"; I = 100103004; J = "Durazno"; K = "American Nectar"; L = "Primera"; M = 220; N = 9500; O = 10000; P = 9773; Q = "$/caja 16 kilos empedrada"; R = "Región de O'Higgins"; S = 611; T = 16 },
    @{ Row = 160; A = 11; B = "Vega Monumental Concepción"; C = "Bíobío"; D = 44595; E = 8; F = "Fruta"; G = 100103; H = "Frutos de hueso (carozo)"; I = 100103004; J = "Durazno"; K = "Elegant Lady"; L = "Primera"; M = 250; N = 10000; O = 11000; P = 10400; Q = "$/caja 16 kilos empedrada"; R = "Región de O'Higgins"; S = 650; T = 16 },
    @{ Row = 161; A = 11; B = "Vega Monumental Concepción"; C = "Bíobío"; D = 44595; E = 8; F = "Fruta"; G = 100103; H = "Frutos de hueso (carozo)"; I = 100103004; J = "Durazno"; K = "Polar King"; L = "Primera"; M = 220; N = 10000; O = 11000; P = 10455; Q = "$/caja 16 kilos empedrada"; R = "Región de O'Higgins"; S = 653; T = 16 }
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    $ws.Cells.Item($rowNum, 1).Value = $r.A
    $ws.Cells.Item($rowNum, 2).Value = $r.B
    $ws.Cells.Item($rowNum, 3).Value = $r.C
    $ws.Cells.Item($rowNum, 4).Value = $r.D
    $ws.Cells.Item($rowNum, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($rowNum, 5).Value = $r.E
    $ws.Cells.Item($rowNum, 6).Value = $r.F
    $ws.Cells.Item($rowNum, 7).Value = $r.G
    $ws.Cells.Item($rowNum, 8).Value = $r.H
    $ws.Cells.Item($rowNum, 9).Value = $r.I
    $ws.Cells.Item($rowNum, 10).Value = $r.J
    $ws.Cells.Item($rowNum, 11).Value = $r.K
    $ws.Cells.Item($rowNum, 12).Value = $r.L
    $ws.Cells.Item($rowNum, 13).Value = $r.M
    $ws.Cells.Item($rowNum, 14).Value = $r.N
    $ws.Cells.Item($rowNum, 15).Value = $r.O
    $ws.Cells.Item($rowNum, 16).Value = $r.P
    $ws.Cells.Item($rowNum, 17).Value = $r.Q
    $ws.Cells.Item($rowNum, 18).Value = $r.R
    $ws.Cells.Item($rowNum, 19).Value = $r.S
    $ws.Cells.Item($rowNum, 20).Value = $r.T
}
